$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New log entries continuing the daily DSA tracker (rows 24-28, 30).
# Shared strings must be introduced in this exact order so new entries land
# on: 47 On Leave, 48 N/A, 49 Row Column Matrix, 50 url, 51 search text.
# ---------------------------------------------------------------------------

# Row 26 - Fri 19 Jan 2024 - On Leave
$ws.Range("A26").Value = 45310
$ws.Range("A23").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B26").Value = "FRI"
$ws.Range("C26").Value = "On Leave"
$ws.Range("D26").Value = "N/A"

# Row 27 - Sat 20 Jan 2024 - On Leave
$ws.Range("A27").Value = 45311
$ws.Range("A23").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B27").Value = "SAT"
$ws.Range("C27").Value = "On Leave"
$ws.Range("D27").Value = "N/A"

# Row 24 - Wed 17 Jan 2024 - Row Column Matrix
$ws.Range("A24").Value = 45308
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B24").Value = "WED"
$ws.Range("C24").Value = "Row Column Matrix"
$ws.Range("D24").Value = "Mediam"

# Row 25 - Thu 18 Jan 2024 - Search element in a sorted matrix
$ws.Range("A25").Value = 45309
$ws.Range("A23").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B25").Value = "THU"
$ws.Range("D25").Value = "Mediam"
$ws.Range("E25").Value = "https://www.geeksforgeeks.org/search-element-sorted-matrix/"
$ws.Hyperlinks.Add($ws.Range("E25"), "https://www.geeksforgeeks.org/search-element-sorted-matrix/")
$ws.Range("E23").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("C25").Value = "Search element in a sorted matrix"

# Row 28 - Sun 21 Jan 2024 - weekly practice / rest day (same style as row 20)
$ws.Range("A28").Value = 45312
$ws.Range("B28").Value = "SUN"
$ws.Range("C28").Value = "Practice All the above"
$ws.Range("A20:E20").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)

# Row 30 - Mon 22 Jan 2024 (row 29 intentionally left blank as a spacer)
$ws.Range("A30").Value = 45313
$ws.Range("A23").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("B30").Value = "MON"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# View state: scroll the window down and select C26, matching where the
# author was working when the file was last saved.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C26").Select()
